# Update pick-rate / win-rate figures on Sheet1 (row 2) with refreshed stats
# after adding all other regions, 20 minute surrenders, and <20 minute games.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 29.91
$ws.Range("B2").Value = 15.64
$ws.Range("C2").Value = 11.56
$ws.Range("D2").Value = 5.86
$ws.Range("E2").Value = 4.25
$ws.Range("F2").Value = 51.01
$ws.Range("G2").Value = 50.35
$ws.Range("H2").Value = 50.97
$ws.Range("I2").Value = 48.43
$ws.Range("J2").Value = 39.67

# Move the active selection as left by the author after editing the sheet.
$ws.Range("I10").Select()
